$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KLAC")

# Row 4 - Inventory
$ws.Range("C4").Value = 1421000000.0
$ws.Range("D4").Value = 1394000000.0
$ws.Range("E4").Value = 1311000000.0
$ws.Range("F4").Value = 1264000000.0
$ws.Range("G4").Value = 1251000000.0

# Row 13 - Accounts Payable
$ws.Range("C13").Value = 262000000.0
$ws.Range("D13").Value = 254000000.0
$ws.Range("E13").Value = 264000000.0
$ws.Range("F13").Value = 243000000.0
$ws.Range("G13").Value = 257000000.0

# Row 19 - Non-current Revenue (Deferred)
$ws.Range("B19").Value = 86902000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("E21").Value = 424000000.0

# Row 22 - Non-current Liabilities (Other)
$ws.Range("B22").Value = 658015000.0

# Row 31 - Shares (Common)
$ws.Range("B31").Value = 153282000.0

# Row 33 - Net Debt
$ws.Range("B33").Value = 999372000.0

# Row 34 - Total Debt
$ws.Range("B34").Value = 3442097000.0
